$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial (45177) for every data row
# (rows 2-295). Update it to 45178 for all of them.
$range = $ws.Range("C2:C295")
$range.Value = 45178
